$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.768.47"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "3.494.13"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.93"
$ws.Range("E5").Value = "  -1.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.18"
$ws.Range("E6").Value = "  +6.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +1.75%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.207"
$ws.Range("E9").Value = "  -2.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.651"
$ws.Range("E10").Value = "  +1.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.99"
$ws.Range("E11").Value = "  +1.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000300"
$ws.Range("E12").Value = "  -1.96%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.51"
$ws.Range("E13").Value = "  +1.09%  "

# Row 14
$ws.Range("D14").Value = "4.049.09"
$ws.Range("E14").Value = "  -0.82%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "607.49"
$ws.Range("E15").Value = "  +5.36%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.890.56"
$ws.Range("E16").Value = "  -0.01%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.08"
$ws.Range("E17").Value = "  +1.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.57"
$ws.Range("E18").Value = "  -0.12%  "

# Row 19
$ws.Range("D19").Value = "3.500.53"
$ws.Range("E19").Value = "  -0.68%  "

# Row 20
$ws.Range("E20").Value = "  +0.61%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.989"
$ws.Range("E21").Value = "  +0.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.10"
$ws.Range("E22").Value = "  +4.80%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.55"
$ws.Range("E23").Value = "  +11.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("E24").Value = "  +4.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.56"
$ws.Range("E25").Value = "  -2.33%  "

# Row 26
$ws.Range("E26").Value = "  +5.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  -0.60%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.71"
$ws.Range("E28").Value = "  +4.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.41"
$ws.Range("E29").Value = "  +4.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.56"
$ws.Range("E30").Value = "  +27.87%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").Value = "  +2.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.64"
$ws.Range("E32").Value = "  +3.95%  "

# Row 33
$ws.Range("E33").Value = "  +1.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.52"
$ws.Range("E34").Value = "  +0.37%  "

# Row 35
$ws.Range("D35").Value = "3.734.80"
$ws.Range("E35").Value = "  +5.61%  "

# Row 36
$ws.Range("E36").Value = "  -0.22%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0803"
$ws.Range("E37").Value = "  +3.29%  "

# Row 38
$ws.Range("E38").Value = "  -10.54%  "

# Row 39
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.390"
$ws.Range("E39").Value = "  -3.13%  "

# Row 40
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.61"
$ws.Range("E40").Value = "  -0.92%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "501.99"
$ws.Range("E41").Value = "  -5.45%  "

# Row 42
$ws.Range("E42").Value = "  +0.71%  "

# Row 43
$ws.Range("E43").Value = "  +0.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0456"
$ws.Range("E44").Value = "  +0.76%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -2.82%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("E46").Value = "  -0.26%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").Value = "  -2.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.71"
$ws.Range("E49").Value = "  -4.07%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.90"
$ws.Range("E50").Value = "  -1.87%  "

# Row 51
$ws.Range("E51").Value = "  +0.59%  "
